$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: expand the "Carl Otto Ehrenfried Nicolai is born..." sentence
# with the birth address, the Kingdom of Prussia, and the mother's father's
# occupation.
# ---------------------------------------------------------------------------
$old1 = "  Carl Otto Ehrenfried Nicolai is born in Königsberg (Kaliningrad), only child of the union of Carl Ernst Daniel Nicolai, a composer, and Christiane Wilhelmine Lauber.  The marriage of his parents will end in a few months owing to the physical and mental condition of his mother.  He will grow up with foster parents until age 10."
$new1 = "  Carl Otto Ehrenfried Nicolai is born at Steindamm 277 in Königsberg (Kaliningrad), Kingdom of Prussia, the only child of the union of Carl Ernst Daniel Nicolai, a composer, and Christiane Wilhelmine Lauber, the daughter of a minister.  The marriage of his parents will end in a few months owing to the physical and mental condition of his mother.  He will grow up with foster parents until age 10."

$found1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    Write-Host "WARNING: change 1 text not found"
}

# ---------------------------------------------------------------------------
# Change 2: add a new "21 June 1810" paragraph (baptism of Otto Nicolai)
# right after the paragraph about the Meyerbeer cantata / Vogler's birthday.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$found2 = $anchor.Find.Execute("their teacher, Georg Joseph Vogler (61).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()

    $newParaStart = $anchor.End + 1

    $dateRange = $d.Range($newParaStart, $newParaStart)
    $dateRange.InsertAfter("21 June 1810")

    $bodyRange = $d.Range($dateRange.End, $dateRange.End)
    $bodyRange.InsertAfter("  Otto Nicolai (0) is baptized in the Lutheran Steindamm Church, Königsberg.")

    $boldRange = $d.Range($newParaStart, $newParaStart + 12)
    $boldRange.Bold = 1
} else {
    Write-Host "WARNING: change 2 anchor text not found"
}

# ---------------------------------------------------------------------------
# Change 3: fix the closing date, "January 2016" -> "May 2016".
# ---------------------------------------------------------------------------
$found3 = $d.Content.Find.Execute("6 January 2016", $false, $false, $false, $false, $false, $true, 1, $false, "6 May 2016", 2)
if (-not $found3) {
    Write-Host "WARNING: change 3 text not found"
}
